$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.04970888862596
$ws.Range("D2").Value = 1.05336706752593
$ws.Range("E2").Value = 1.04680113597715
$ws.Range("F2").Value = 1.061627691383946
$ws.Range("I2").Value = 1.038302735358962
$ws.Range("J2").Value = 1.054745901642773
$ws.Range("K2").Value = 1.056113232328339
$ws.Range("L2").Value = 1.049565555426539
$ws.Range("M2").Value = 1.064351237256879

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.051408124854707
$ws.Range("D3").Value = 1.054710776521811
$ws.Range("E3").Value = 1.048273350027508
$ws.Range("F3").Value = 1.063200686143314
$ws.Range("I3").Value = 1.038686550608598
$ws.Range("J3").Value = 1.056090861217224
$ws.Range("K3").Value = 1.057268612926101
$ws.Range("L3").Value = 1.050847771178331
$ws.Range("M3").Value = 1.065736980932223

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.052505325776798
$ws.Range("D4").Value = 1.055577957420256
$ws.Range("E4").Value = 1.049224047233271
$ws.Range("F4").Value = 1.064216531138523
$ws.Range("I4").Value = 1.038932505001578
$ws.Range("J4").Value = 1.056958514722631
$ws.Range("K4").Value = 1.058013399260951
$ws.Range("L4").Value = 1.051675023326691
$ws.Range("M4").Value = 1.066631172740828

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.052966046494882
$ws.Range("D5").Value = 1.055941981271522
$ws.Range("E5").Value = 1.049623269633533
$ws.Range("F5").Value = 1.0646431269236
$ws.Range("I5").Value = 1.039035332959586
$ws.Range("J5").Value = 1.057322657464132
$ws.Range("K5").Value = 1.058325840118522
$ws.Range("L5").Value = 1.052022228145743
$ws.Range("M5").Value = 1.067006507577722

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.053043372077187
$ws.Range("D6").Value = 1.056003071128437
$ws.Range("E6").Value = 1.049690274686087
$ws.Range("F6").Value = 1.064714727240327
$ws.Range("I6").Value = 1.039052564820105
$ws.Range("J6").Value = 1.057383762626948
$ws.Range("K6").Value = 1.05837826142179
$ws.Range("L6").Value = 1.052080492083908
$ws.Range("M6").Value = 1.067069494015104

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.05251148406419
$ws.Range("D7").Value = 1.055582823630303
$ws.Range("E7").Value = 1.049229383419869
$ws.Range("F7").Value = 1.064222233146748
$ws.Range("I7").Value = 1.038933881233966
$ws.Range("J7").Value = 1.056963382833779
$ws.Range("K7").Value = 1.058017576721558
$ws.Range("L7").Value = 1.051679664929409
$ws.Range("M7").Value = 1.066636190262355

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.050283641094671
$ws.Range("D8").Value = 1.053821659544518
$ws.Range("E8").Value = 1.047299081874205
$ws.Range("F8").Value = 1.062159710217082
$ws.Range("I8").Value = 1.038432946356402
$ws.Range("J8").Value = 1.05520098647182
$ws.Range("K8").Value = 1.056504287985327
$ws.Range("L8").Value = 1.049999394298822
$ws.Range("M8").Value = 1.064820074503733

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.046339576494526
$ws.Range("D9").Value = 1.050700333107755
$ws.Range("E9").Value = 1.043882467239555
$ws.Range("F9").Value = 1.058509587237848
$ws.Range("I9").Value = 1.03753170998789
$ws.Range("J9").Value = 1.052074875273024
$ws.Range("K9").Value = 1.053815690597578
$ws.Range("L9").Value = 1.047019543998641
$ws.Range("M9").Value = 1.061600442908638

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.043697091475469
$ws.Range("D10").Value = 1.048606828285334
$ws.Range("E10").Value = 1.041593897322076
$ws.Range("F10").Value = 1.056064959705235
$ws.Range("I10").Value = 1.036918224592825
$ws.Range("J10").Value = 1.049976374322297
$ws.Range("K10").Value = 1.052007980689353
$ws.Range("L10").Value = 1.045019627900092
$ws.Range("M10").Value = 1.059440367722662

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.042549580245486
$ws.Range("D11").Value = 1.047697197409005
$ws.Range("E11").Value = 1.040600213796262
$ws.Range("F11").Value = 1.055003601685753
$ws.Range("I11").Value = 1.036649529227562
$ws.Range("J11").Value = 1.049064140758654
$ws.Range("K11").Value = 1.051221472363606
$ws.Range("L11").Value = 1.044150346952669
$ws.Range("M11").Value = 1.058501656444982

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.042122832705172
$ws.Range("D12").Value = 1.047358838809663
$ws.Range("E12").Value = 1.040230695085915
$ws.Range("F12").Value = 1.054608929823063
$ws.Range("I12").Value = 1.03654926119009
$ws.Range("J12").Value = 1.04872474846828
$ws.Range("K12").Value = 1.05092875326277
$ws.Range("L12").Value = 1.043826949702147
$ws.Range("M12").Value = 1.058152456959303

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.042214394862572
$ws.Range("D13").Value = 1.047431439830911
$ws.Range("E13").Value = 1.04030997732361
$ws.Range("F13").Value = 1.054693608225574
$ws.Range("I13").Value = 1.036570790024559
$ws.Range("J13").Value = 1.048797574263854
$ws.Range("K13").Value = 1.050991568684905
$ws.Range("L13").Value = 1.043896342661315
$ws.Range("M13").Value = 1.058227385171448

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.042514315669563
$ws.Range("D14").Value = 1.047669238461959
$ws.Range("E14").Value = 1.040569677911027
$ws.Range("F14").Value = 1.054970986960327
$ws.Range("I14").Value = 1.036641250501196
$ws.Range("J14").Value = 1.049036097737301
$ws.Range("K14").Value = 1.051197287898914
$ws.Range("L14").Value = 1.044123625255199
$ws.Range("M14").Value = 1.058472802189909

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.042699038633473
$ws.Range("D15").Value = 1.047815689952705
$ws.Range("E15").Value = 1.040729631941974
$ws.Range("F15").Value = 1.055141830961385
$ws.Range("I15").Value = 1.036684602102939
$ws.Range("J15").Value = 1.049182986929536
$ws.Range("K15").Value = 1.051323961808072
$ws.Range("L15").Value = 1.044263593963279
$ws.Range("M15").Value = 1.058623942410754

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.043773177231834
$ws.Range("D16").Value = 1.048667130618234
$ws.Range("E16").Value = 1.041659786536609
$ws.Range("F16").Value = 1.056135338073073
$ws.Range("I16").Value = 1.036935992368214
$ws.Range("J16").Value = 1.050036840051349
$ws.Range("K16").Value = 1.052060098603827
$ws.Range("L16").Value = 1.045077248647074
$ws.Range("M16").Value = 1.059502594583707

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.044446062318362
$ws.Range("D17").Value = 1.049200371259039
$ws.Range("E17").Value = 1.042242511838027
$ws.Range("F17").Value = 1.056757775310128
$ws.Range("I17").Value = 1.037092862839225
$ws.Range("J17").Value = 1.050571475684803
$ws.Range("K17").Value = 1.052520844411018
$ws.Range("L17").Value = 1.045586740266392
$ws.Range("M17").Value = 1.06005283543779

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.044838227660145
$ws.Range("D18").Value = 1.049511100305223
$ws.Range("E18").Value = 1.042582144258262
$ws.Range("F18").Value = 1.057120561381946
$ws.Range("I18").Value = 1.037184068520771
$ws.Range("J18").Value = 1.050882976080774
$ws.Range("K18").Value = 1.052789227753086
$ws.Range("L18").Value = 1.045883600075098
$ws.Range("M18").Value = 1.060373455873155

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.04497189257511
$ws.Range("D19").Value = 1.049617000123166
$ws.Range("E19").Value = 1.042697906164855
$ws.Range("F19").Value = 1.057244216470794
$ws.Range("I19").Value = 1.037215117520278
$ws.Range("J19").Value = 1.050989131758015
$ws.Range("K19").Value = 1.052880678463466
$ws.Range("L19").Value = 1.045984767991419
$ws.Range("M19").Value = 1.060482724260236

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.044373901031186
$ws.Range("D20").Value = 1.049143190818797
$ws.Range("E20").Value = 1.042180018052263
$ws.Range("F20").Value = 1.056691021812783
$ws.Range("I20").Value = 1.037076062583216
$ws.Range("J20").Value = 1.050514149959734
$ws.Range("K20").Value = 1.05247144825293
$ws.Range("L20").Value = 1.045532109614007
$ws.Range("M20").Value = 1.059993833580548

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.042426010690991
$ws.Range("D21").Value = 1.047599226064408
$ws.Range("E21").Value = 1.040493214264006
$ws.Range("F21").Value = 1.054889317982485
$ws.Range("I21").Value = 1.036620514450324
$ws.Range("J21").Value = 1.048965873747045
$ws.Range("K21").Value = 1.051136724670756
$ws.Range("L21").Value = 1.044056710299001
$ws.Range("M21").Value = 1.058400547453962

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.041198331824136
$ws.Range("D22").Value = 1.046625684578137
$ws.Range("E22").Value = 1.03943021534984
$ws.Range("F22").Value = 1.053753984883722
$ws.Range("I22").Value = 1.036331414663116
$ws.Range("J22").Value = 1.047989233731397
$ws.Range("K22").Value = 1.050294199669673
$ws.Range("L22").Value = 1.043126125729925
$ws.Range("M22").Value = 1.057395768855433

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.041849433623997
$ws.Range("D23").Value = 1.047142045468665
$ws.Range("E23").Value = 1.039993966396083
$ws.Range("F23").Value = 1.054356090613272
$ws.Range("I23").Value = 1.036484927199071
$ws.Range("J23").Value = 1.048507274277905
$ws.Range("K23").Value = 1.050741157347508
$ws.Range("L23").Value = 1.043619728898149
$ws.Range("M23").Value = 1.057928710904593

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.044406508590584
$ws.Range("D24").Value = 1.049169029129156
$ws.Range("E24").Value = 1.042208257111721
$ws.Range("F24").Value = 1.056721185682366
$ws.Range("I24").Value = 1.037083654804856
$ws.Range("J24").Value = 1.050540054048781
$ws.Range("K24").Value = 1.052493769368766
$ws.Range("L24").Value = 1.045556795836963
$ws.Range("M24").Value = 1.060020494986668

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.047361463783019
$ws.Range("D25").Value = 1.051509453478216
$ws.Range("E25").Value = 1.044767606286392
$ws.Range("F25").Value = 1.059455158941107
$ws.Range("I25").Value = 1.037766917571453
$ws.Range("J25").Value = 1.052885547361145
$ws.Range("K25").Value = 1.054513417266078
$ws.Range("L25").Value = 1.047792217082303
$ws.Range("M25").Value = 1.062435156523322
